$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1325
$ws1.Range("F4").Value = 14523
$ws1.Range("F5").Value = 17379
$ws1.Range("F10").Value = 207
$ws1.Range("F12").Value = 53
$ws1.Range("F16").Value = 40
$ws1.Range("F25").Value = 7142
$ws1.Range("F28").Value = 1161
$ws1.Range("F30").Value = 5837
$ws1.Range("F31").Value = 57
$ws1.Range("F32").Value = 36
$ws1.Range("F36").Value = 5036

# Sheet "全部类型" (sheet4) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1325
$ws4.Range("F4").Value = 14523
$ws4.Range("F5").Value = 17379
$ws4.Range("F10").Value = 207
$ws4.Range("F12").Value = 53
$ws4.Range("F16").Value = 40
$ws4.Range("F26").Value = 7142
$ws4.Range("F29").Value = 1161
$ws4.Range("F32").Value = 5837
$ws4.Range("F33").Value = 57
$ws4.Range("F34").Value = 36
$ws4.Range("F38").Value = 5036
